# Update the NroSiniestro (claim numbers) test data for the cancellation
# (anulación) test case: two new smart folders are added, one to perform
# the cancellation itself and another one to obtain the cancellation
# number, so the existing rows get fresh claim numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these numeric-looking claim numbers to stay
# text (same as how they were already stored), preserving leading zeros
# and the trailing space on the "Juicio" row's claim number.
$ws.Range("F4").Value = "'0420194406824"
$ws.Range("F3").Value = "'0420172008616 "
$ws.Range("F2").Value = "'0420172008636"

$ws.Range("G8").Select()
